$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 23:29"

# --- Country reordering: update country names (column A) for rows whose
#     position in the country list shifted ---
$ws.Range("A126").Value = "Ruanda"  # was "Mali"
$ws.Range("A127").Value = "Mali"  # was "Ruanda"
$ws.Range("A128").Value = "Tunez"  # was "Eslovenia"
$ws.Range("A129").Value = "Eslovenia"  # was "Lituania"
$ws.Range("A130").Value = "Lituania"  # was "Sudan del Sur"
$ws.Range("A131").Value = "Sudan del Sur"  # was "Tunez"
$ws.Range("A147").Value = "Republica de Chipre"  # was "Aruba"
$ws.Range("A148").Value = "Aruba"  # was "Republica de Chipre"
$ws.Range("A154").Value = "Togo"  # was "Jamaica"
$ws.Range("A155").Value = "Jamaica"  # was "Togo"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 5737918
$ws.Range("C4").Value = 36987
$ws.Range("D4").Value = 3080189
$ws.Range("E4").Value = 2480551
$ws.Range("G4").Value = 844
$ws.Range("H4").Value = 177178

$ws.Range("B5").Value = 3501975
$ws.Range("C5").Value = 41562
$ws.Range("E5").Value = 774417
$ws.Range("G5").Value = 1115
$ws.Range("H5").Value = 112304

$ws.Range("B8").Value = 599940
$ws.Range("C8").Value = 3880
$ws.Range("D8").Value = 497169
$ws.Range("E8").Value = 90153
$ws.Range("G8").Value = 195
$ws.Range("H8").Value = 12618

$ws.Range("B22").Value = 231284
$ws.Range("C22").Value = 1584
$ws.Range("E22").Value = 17160

$ws.Range("B32").Value = 99599
$ws.Range("C32").Value = 1630
$ws.Range("D32").Value = 74579
$ws.Range("E32").Value = 24225

$ws.Range("B77").Value = 17249
$ws.Range("C77").Value = 17
$ws.Range("D77").Value = 14611
$ws.Range("E77").Value = 2526
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 112

$ws.Range("B93").Value = 8876
$ws.Range("C93").Value = 84
$ws.Range("D93").Value = 7628
$ws.Range("E93").Value = 1195

$ws.Range("B126").Value = 2717
$ws.Range("C126").Value = 73
$ws.Range("D126").Value = 1705
$ws.Range("E126").Value = 1001
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 11

$ws.Range("B127").Value = 2667
$ws.Range("D127").Value = 1993
$ws.Range("E127").Value = 549
$ws.Range("H127").Value = 125

$ws.Range("B128").Value = 2543
$ws.Range("C128").Value = 116
$ws.Range("D128").Value = 1397
$ws.Range("E128").Value = 1083
$ws.Range("G128").Value = 3
$ws.Range("H128").Value = 63

$ws.Range("B129").Value = 2536
$ws.Range("C129").Value = 43
$ws.Range("D129").Value = 2079
$ws.Range("E129").Value = 328
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 129

$ws.Range("B130").Value = 2528
$ws.Range("C130").Value = 32
$ws.Range("D130").Value = 1747
$ws.Range("E130").Value = 699
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 82

$ws.Range("B131").Value = 2494
$ws.Range("D131").Value = 1290
$ws.Range("E131").Value = 1157
$ws.Range("H131").Value = 47

$ws.Range("B147").Value = 1395
$ws.Range("C147").Value = 10
$ws.Range("D147").Value = 878
$ws.Range("E147").Value = 497
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 20

$ws.Range("B148").Value = 1387
$ws.Range("C148").Value = 91
$ws.Range("D148").Value = 277
$ws.Range("E148").Value = 1104
$ws.Range("G148").Value = 1
$ws.Range("H148").Value = 6

$ws.Range("B154").Value = 1212
$ws.Range("C154").Value = 22
$ws.Range("D154").Value = 878
$ws.Range("E154").Value = 307
$ws.Range("H154").Value = 27

$ws.Range("B155").Value = 1192
$ws.Range("C155").Value = 46
$ws.Range("D155").Value = 772
$ws.Range("E155").Value = 406
$ws.Range("H155").Value = 14

$ws.Range("D156").Value = 1082
$ws.Range("E156").Value = 16

$ws.Range("B162").Value = 888
$ws.Range("C162").Value = 3
$ws.Range("E162").Value = 44

